$wb = $excel.ActiveWorkbook


# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 74
$ws.Range("H74").Value = 2999349.2
$ws.Range("I74").Value = 3398662.5
$ws.Range("J74").Value = 4500
$ws.Range("K74").Value = 3398662.5
$ws.Range("L74").Value = 4500
$ws.Range("M74").Value = -3397726.5
$ws.Range("N74").Value = -6372
# Row 77
$ws.Range("H77").Value = 2999349.2
$ws.Range("I77").Value = 3398662.5
$ws.Range("J77").Value = 4500
$ws.Range("K77").Value = 16993312.5
$ws.Range("L77").Value = 22500
$ws.Range("M77").Value = -16988632.5
$ws.Range("N77").Value = -31860
# Row 139
$ws.Range("H139").Value = 52170
$ws.Range("J139").Value = 52170
$ws.Range("L139").Value = 52170
$ws.Range("N139").Value = -62450

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 9219.036
$ws.Range("I32").Value = 4143.8823
$ws.Range("J32").Value = 60985.6
$ws.Range("K32").Value = 4143.8823
$ws.Range("L32").Value = 60985.6
$ws.Range("M32").Value = -3856.8823
$ws.Range("N32").Value = -61559.6
# Row 49
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
# Row 61
$ws.Range("H61").Value = 1390
$ws.Range("I61").Value = 1012
$ws.Range("J61").Value = 1957
$ws.Range("K61").Value = 1012
$ws.Range("L61").Value = 1957
$ws.Range("M61").Value = -800
$ws.Range("N61").Value = -2381
# Row 74
$ws.Range("H74").Value = 796.325
$ws.Range("I74").Value = 734.5185
$ws.Range("K74").Value = 734.5185
$ws.Range("M74").Value = 139.4815
# Row 77
$ws.Range("H77").Value = 796.325
$ws.Range("I77").Value = 734.5185
$ws.Range("K77").Value = 3672.5925
$ws.Range("M77").Value = 695.4074999999998
# Row 82
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
# Row 85
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
# Row 136
$ws.Range("H136").Value = 1390
$ws.Range("I136").Value = 1012
$ws.Range("J136").Value = 1957
$ws.Range("K136").Value = 3036
$ws.Range("L136").Value = 5871
$ws.Range("M136").Value = -486
$ws.Range("N136").Value = -10971

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 135
$ws.Range("H135").Value = 64264.75
$ws.Range("J135").Value = 64264.75
$ws.Range("L135").Value = 64264.75
$ws.Range("N135").Value = -74404.75
# Row 137
$ws.Range("H137").Value = 51175
$ws.Range("J137").Value = 51175
$ws.Range("L137").Value = 51175
$ws.Range("N137").Value = -61375
# Row 138
$ws.Range("H138").Value = 54589.09
$ws.Range("J138").Value = 54589.09
$ws.Range("L138").Value = 54589.09
$ws.Range("N138").Value = -64869.09

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2460.4783
$ws.Range("I31").Value = 2403.65
$ws.Range("J31").Value = 2839.3333
$ws.Range("K31").Value = 2403.65
$ws.Range("L31").Value = 2839.3333
$ws.Range("M31").Value = -2108.65
$ws.Range("N31").Value = -3429.3333
# Row 33
$ws.Range("H33").Value = 1311.375
$ws.Range("I33").Value = 927.2857
$ws.Range("J33").Value = 4000
$ws.Range("K33").Value = 927.2857
$ws.Range("L33").Value = 4000
$ws.Range("M33").Value = -548.2857
$ws.Range("N33").Value = -4758
# Row 34
$ws.Range("H34").Value = 2460.4783
$ws.Range("I34").Value = 2403.65
$ws.Range("J34").Value = 2839.3333
$ws.Range("K34").Value = 2403.65
$ws.Range("L34").Value = 2839.3333
$ws.Range("M34").Value = -2201.65
$ws.Range("N34").Value = -3243.3333
# Row 58
$ws.Range("H58").Value = 5049
$ws.Range("I58").Value = 790.26086
$ws.Range("J58").Value = 103000
$ws.Range("K58").Value = 790.26086
$ws.Range("L58").Value = 103000
$ws.Range("M58").Value = -587.26086
$ws.Range("N58").Value = -103406
# Row 122
$ws.Range("H122").Value = 3381.5
$ws.Range("I122").Value = 5506
$ws.Range("J122").Value = 1257
$ws.Range("K122").Value = 16518
$ws.Range("L122").Value = 3771
$ws.Range("M122").Value = -14068
$ws.Range("N122").Value = -8671
# Row 136
$ws.Range("H136").Value = 5049
$ws.Range("I136").Value = 790.26086
$ws.Range("J136").Value = 103000
$ws.Range("K136").Value = 2370.78258
$ws.Range("L136").Value = 309000
$ws.Range("M136").Value = 179.2174199999999
$ws.Range("N136").Value = -314100

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 58
$ws.Range("H58").Value = 3381.818
$ws.Range("I58").Value = 500
$ws.Range("J58").Value = 4022.2222
$ws.Range("K58").Value = 1500
$ws.Range("L58").Value = 12066.6666
$ws.Range("M58").Value = -1372
$ws.Range("N58").Value = -12322.6666

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 65
$ws.Range("I2").Value = 49
$ws.Range("J2").Value = 81
$ws.Range("K2").Value = 49
$ws.Range("L2").Value = 81
$ws.Range("M2").Value = 64
$ws.Range("N2").Value = -307
# Row 36
$ws.Range("H36").Value = 2964.2856
$ws.Range("J36").Value = 2961.5386
$ws.Range("L36").Value = 2961.5386
$ws.Range("N36").Value = -3931.5386
# Row 43
$ws.Range("H43").Value = 9933.333000000001
$ws.Range("I43").Value = 7500
$ws.Range("K43").Value = 7500
$ws.Range("M43").Value = -7349
# Row 46
$ws.Range("H46").Value = 15960
$ws.Range("J46").Value = 22600
$ws.Range("L46").Value = 22600
$ws.Range("N46").Value = -22912
# Row 70
$ws.Range("H70").Value = 5368.3076
$ws.Range("I70").Value = 4866.6665
$ws.Range("K70").Value = 4866.6665
$ws.Range("M70").Value = -4596.6665
# Row 73
$ws.Range("H73").Value = 5368.3076
$ws.Range("I73").Value = 4866.6665
$ws.Range("K73").Value = 4866.6665
$ws.Range("M73").Value = -3930.6665
# Row 136
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
# Row 138
$ws.Range("H138").Value = 43784.832
$ws.Range("J138").Value = 45541.8
$ws.Range("L138").Value = 45541.8
$ws.Range("N138").Value = -55821.8

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 1510
$ws.Range("I46").Value = 2325
$ws.Range("J46").Value = 695
$ws.Range("K46").Value = 2325
$ws.Range("L46").Value = 695
$ws.Range("M46").Value = -2137
$ws.Range("N46").Value = -1071
# Row 64
$ws.Range("H64").Value = 19333.334
$ws.Range("J64").Value = 19333.334
$ws.Range("L64").Value = 19333.334
$ws.Range("N64").Value = -19783.334
# Row 67
$ws.Range("H67").Value = 19333.334
$ws.Range("J67").Value = 19333.334
$ws.Range("L67").Value = 19333.334
$ws.Range("N67").Value = -20893.334
# Row 82
$ws.Range("H82").Value = 3432.2
$ws.Range("I82").Value = 3244.4443
$ws.Range("J82").Value = 3713.8333
$ws.Range("K82").Value = 3244.4443
$ws.Range("L82").Value = 3713.8333
$ws.Range("M82").Value = -2883.4443
$ws.Range("N82").Value = -4435.8333
# Row 85
$ws.Range("H85").Value = 3432.2
$ws.Range("I85").Value = 3244.4443
$ws.Range("J85").Value = 3713.8333
$ws.Range("K85").Value = 3244.4443
$ws.Range("L85").Value = 3713.8333
$ws.Range("M85").Value = -1996.4443
$ws.Range("N85").Value = -6209.8333
# Row 94
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
# Row 132
$ws.Range("H132").Value = 2384.5217
$ws.Range("I132").Value = 1359.6296
$ws.Range("J132").Value = 3840.9473
$ws.Range("K132").Value = 4078.8888
$ws.Range("L132").Value = 11522.8419
$ws.Range("M132").Value = -1548.8888
$ws.Range("N132").Value = -16582.8419
# Row 138
$ws.Range("H138").Value = 38750
$ws.Range("J138").Value = 38750
$ws.Range("L138").Value = 38750
$ws.Range("N138").Value = -49030

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Range("H54").Value = 100000
$ws.Range("J54").Value = 100000
$ws.Range("L54").Value = 100000
$ws.Range("N54").Value = -101040
# Row 63
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
# Row 66
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
# Row 137
$ws.Range("H137").Value = 65400
$ws.Range("J137").Value = 65400
$ws.Range("L137").Value = 65400
$ws.Range("N137").Value = -75600
# Row 138
$ws.Range("H138").Value = 54196
$ws.Range("J138").Value = 54196
$ws.Range("L138").Value = 54196
$ws.Range("N138").Value = -64476
# Row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
